# Update BAP Yearly Financials sheet with latest reported figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 3326300
$ws.Range("E8").Value = 3248600
$ws.Range("F8").Value = 2950400
$ws.Range("G8").Value = 2593600
$ws.Range("H8").Value = 2136900
$ws.Range("I8").Value = 1836900
$ws.Range("J8").Value = 554200

# Row 15
$ws.Range("D15").Value = -126600
$ws.Range("E15").Value = -122700
$ws.Range("F15").Value = -120900
$ws.Range("G15").Value = -130800
$ws.Range("H15").Value = -135500
$ws.Range("I15").Value = -119000
$ws.Range("J15").Value = -28300

# Row 17
$ws.Range("D17").Value = 1431900
$ws.Range("E17").Value = 1417300
$ws.Range("F17").Value = 1329200
$ws.Range("G17").Value = 1178100
$ws.Range("H17").Value = 1009300
$ws.Range("I17").Value = 851900
$ws.Range("J17").Value = 225100

# Row 18
$ws.Range("D18").Value = 1894400
$ws.Range("E18").Value = 1831300
$ws.Range("F18").Value = 1621200
$ws.Range("G18").Value = 1415500
$ws.Range("H18").Value = 1127700
$ws.Range("I18").Value = 985000
$ws.Range("J18").Value = 329100

# Row 20
$ws.Range("D20").Value = -213300
$ws.Range("E20").Value = -356300
$ws.Range("F20").Value = -306200
$ws.Range("G20").Value = -393400
$ws.Range("H20").Value = -424200
$ws.Range("I20").Value = -142400
$ws.Range("J20").Value = -47300

# Row 21
$ws.Range("D21").Value = 1809700
$ws.Range("E21").Value = 1597800
$ws.Range("F21").Value = 1434500
$ws.Range("G21").Value = 1152900
$ws.Range("H21").Value = 802500
$ws.Range("I21").Value = 875300
$ws.Range("J21").Value = "NA"

# Row 23
$ws.Range("D23").Value = 1681100
$ws.Range("E23").Value = 1475000
$ws.Range("F23").Value = 1314900
$ws.Range("G23").Value = 1022100
$ws.Range("H23").Value = 703400
$ws.Range("I23").Value = 842600
$ws.Range("J23").Value = 281700

# Row 24
$ws.Range("D24").Value = 420100
$ws.Range("E24").Value = 386400
$ws.Range("F24").Value = 361000
$ws.Range("G24").Value = 292000
$ws.Range("H24").Value = 233800
$ws.Range("I24").Value = 200000
$ws.Range("J24").Value = 63500

# Row 26
$ws.Range("D26").Value = 1261000
$ws.Range("E26").Value = 1088600
$ws.Range("F26").Value = 953900
$ws.Range("G26").Value = 730100
$ws.Range("H26").Value = 469700
$ws.Range("I26").Value = 642600
$ws.Range("J26").Value = 218300

# Row 27
$ws.Range("D27").Value = 1233900
$ws.Range("E27").Value = 1059800
$ws.Range("F27").Value = 932500
$ws.Range("G27").Value = 720100
$ws.Range("H27").Value = 463900
$ws.Range("I27").Value = 627100
$ws.Range("J27").Value = 213900

# Row 32
$ws.Range("D32").Value = 213300
$ws.Range("E32").Value = 356300
$ws.Range("F32").Value = 306200
$ws.Range("G32").Value = 393400
$ws.Range("H32").Value = 424200
$ws.Range("I32").Value = 142400
$ws.Range("J32").Value = 47300

# Row 33
$ws.Range("D33").Value = 1233900
$ws.Range("E33").Value = 1059800
$ws.Range("F33").Value = 932500
$ws.Range("G33").Value = 720100
$ws.Range("H33").Value = 463900
$ws.Range("I33").Value = 627100
$ws.Range("J33").Value = 213900

# Row 35
$ws.Range("D35").Value = 1233900
$ws.Range("E35").Value = 1059800
$ws.Range("F35").Value = 932500
$ws.Range("G35").Value = 720100
$ws.Range("H35").Value = 463900
$ws.Range("I35").Value = 627100
$ws.Range("J35").Value = 213900

# Row 41
$ws.Range("D41").Value = 7002600
$ws.Range("E41").Value = 5019500
$ws.Range("F41").Value = 6752200
$ws.Range("G41").Value = 14233800
$ws.Range("H41").Value = 6562600
$ws.Range("I41").Value = 6381500
$ws.Range("J41").Value = 330100

# Row 42
$ws.Range("D42").Value = 3681000
$ws.Range("E42").Value = 4787800
$ws.Range("F42").Value = 4470600
$ws.Range("G42").Value = 2778000
$ws.Range("H42").Value = 1068700
$ws.Range("I42").Value = 3651000
$ws.Range("J42").Value = 1356500

# Row 47
$ws.Range("D47").Value = 213800
$ws.Range("E47").Value = 211000
$ws.Range("F47").Value = 190200
$ws.Range("G47").Value = 49000
$ws.Range("H47").Value = 19200
$ws.Range("I47").Value = 17900
$ws.Range("J47").Value = 4000

# Row 48
$ws.Range("D48").Value = 593600
$ws.Range("E48").Value = 620200
$ws.Range("F48").Value = 631000
$ws.Range("G48").Value = 1260100
$ws.Range("H48").Value = 1423000
$ws.Range("I48").Value = 785900
$ws.Range("J48").Value = 142500

# Row 49
$ws.Range("D49").Value = 596700
$ws.Range("E49").Value = 591200
$ws.Range("F49").Value = 573200
$ws.Range("G49").Value = 998400
$ws.Range("H49").Value = 748600
$ws.Range("I49").Value = 813200
$ws.Range("J49").Value = 136700

# Row 52
$ws.Range("D52").Value = 144800
$ws.Range("E52").Value = 152600
$ws.Range("F52").Value = 158200
$ws.Range("G52").Value = 284800
$ws.Range("H52").Value = 102700
$ws.Range("I52").Value = 87200
$ws.Range("J52").Value = 29000

# Row 54
$ws.Range("D54").Value = 51405900
$ws.Range("E54").Value = 47173000
$ws.Range("F54").Value = 46885100
$ws.Range("G54").Value = 40659300
$ws.Range("H54").Value = 34405100
$ws.Range("I54").Value = 31371000
$ws.Range("J54").Value = 9261900

# Row 57
$ws.Range("D57").Value = 636500
$ws.Range("E57").Value = 524900
$ws.Range("F57").Value = 465200
$ws.Range("G57").Value = 5989800
$ws.Range("H57").Value = 1382000
$ws.Range("I57").Value = 1579900
$ws.Range("J57").Value = 22700

# Row 59
$ws.Range("D59").Value = 339200
$ws.Range("E59").Value = 292000
$ws.Range("F59").Value = 297900
$ws.Range("G59").Value = 228600
$ws.Range("H59").Value = 159300
$ws.Range("I59").Value = 121300
$ws.Range("J59").Value = 46800

# Row 61
$ws.Range("D61").Value = 7238300
$ws.Range("E61").Value = 6997700
$ws.Range("F61").Value = 7194000
$ws.Range("G61").Value = 7334300
$ws.Range("H61").Value = 4262000
$ws.Range("I61").Value = 3678200
$ws.Range("J61").Value = 1195800

# Row 62
$ws.Range("D62").Value = 261900
$ws.Range("E62").Value = 204600
$ws.Range("F62").Value = 158900
$ws.Range("G62").Value = 1952300
$ws.Range("H62").Value = 1679000
$ws.Range("I62").Value = 625900
$ws.Range("J62").Value = 457400

# Row 66
$ws.Range("D66").Value = 44845200
$ws.Range("E66").Value = 41245700
$ws.Range("F66").Value = 42021700
$ws.Range("G66").Value = 36443800
$ws.Range("H66").Value = 30837300
$ws.Range("I66").Value = 28166100
$ws.Range("J66").Value = 8237900

# Row 72
$ws.Range("D72").Value = 5705000
$ws.Range("E72").Value = 5143200
$ws.Range("F72").Value = 4213000
$ws.Range("G72").Value = 3789200
$ws.Range("H72").Value = 3149700
$ws.Range("I72").Value = 3765600
$ws.Range("J72").Value = 768200

# Row 76
$ws.Range("D76").Value = 6560700
$ws.Range("E76").Value = 5927300
$ws.Range("F76").Value = 4863400
$ws.Range("G76").Value = 4215500
$ws.Range("H76").Value = 3567800
$ws.Range("I76").Value = 3205000
$ws.Range("J76").Value = 1024000

# Row 81
$ws.Range("D81").Value = 1233900
$ws.Range("E81").Value = 1059800
$ws.Range("F81").Value = 932500
$ws.Range("G81").Value = 720100
$ws.Range("H81").Value = 463900
$ws.Range("I81").Value = 627100
$ws.Range("J81").Value = 213900

# Row 83
$ws.Range("D83").Value = 128600
$ws.Range("E83").Value = 122700
$ws.Range("F83").Value = 119600
$ws.Range("G83").Value = 130800
$ws.Range("H83").Value = 99000
$ws.Range("I83").Value = 32700
$ws.Range("J83").Value = "NA"

# Row 89
$ws.Range("D89").Value = 2676600
$ws.Range("E89").Value = -439500
$ws.Range("F89").Value = -416700
$ws.Range("G89").Value = -550400
$ws.Range("H89").Value = 1104900
$ws.Range("I89").Value = 1081100
$ws.Range("J89").Value = -497600

# Row 91
$ws.Range("D91").Value = -43400
$ws.Range("E91").Value = -33200
$ws.Range("F91").Value = -44800
$ws.Range("G91").Value = -91000
$ws.Range("H91").Value = -162200
$ws.Range("I91").Value = -48900
$ws.Range("J91").Value = -30400

# Row 94
$ws.Range("D94").Value = 122300
$ws.Range("E94").Value = -569600
$ws.Range("F94").Value = -457200
$ws.Range("G94").Value = 693400
$ws.Range("H94").Value = -557100
$ws.Range("I94").Value = -645600
$ws.Range("J94").Value = "NA"

# Row 96
$ws.Range("D96").Value = -673100
$ws.Range("E96").Value = -197000
$ws.Range("F96").Value = -162800
$ws.Range("G96").Value = -129500
$ws.Range("H96").Value = -161400
$ws.Range("I96").Value = -55300
$ws.Range("J96").Value = -46900

# Row 100
$ws.Range("D100").Value = -706600
$ws.Range("E100").Value = -485800
$ws.Range("F100").Value = 16200
$ws.Range("G100").Value = -302500
$ws.Range("H100").Value = 201600
$ws.Range("I100").Value = 237800
$ws.Range("J100").Value = "NA"

# Row 101
$ws.Range("D101").Value = -236600
$ws.Range("E101").Value = -136900
$ws.Range("F101").Value = 964900
$ws.Range("G101").Value = 137300
$ws.Range("H101").Value = -222400
$ws.Range("I101").Value = 34200
$ws.Range("J101").Value = "NA"

# Row 102
$ws.Range("D102").Value = 1983900
$ws.Range("E102").Value = -1631900
$ws.Range("F102").Value = 107200
$ws.Range("G102").Value = -22200
$ws.Range("H102").Value = 527000
$ws.Range("I102").Value = 707500
$ws.Range("J102").Value = -917200
